# This script updates the Customer_Reviews column (C2:C51) so that the
# previously-stored full-precision decimal values are replaced with the
# values rounded to the nearest whole number, and removes the thin border
# that was applied around the header row (A1:D1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rounded Customer_Reviews values (column C), rows 2-51.
$roundedValues = @(
    222, 503, 774, 1005, 1287, 1463, 1765, 1981, 2225, 2485,
    821, 2114, 1945, 740, 2197, 2123, 1921, 880, 810, 1792,
    825, 1681, 564, 1903, 1788, 2299, 1979, 1127, 1540, 958,
    1594, 1603, 1539, 2023, 1424, 2344, 2246, 940, 625, 1670,
    755, 1369, 2118, 1263, 1771, 1005, 621, 1199, 1584, 1667
)

$startRow = 2
for ($i = 0; $i -lt $roundedValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $roundedValues[$i]
}

# Remove the thin border that outlined the header cells.
$ws.Range("A1:D1").Borders.LineStyle = 0
